# Apply weekly fruit/vegetable price updates to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date and volume update
$ws.Range("D2").Value = 44306
$ws.Range("M2").Value = 80

# Row 3: date, volume, and commercialization unit/price details update
$ws.Range("D3").Value = 44309
$ws.Range("M3").Value = 80
$ws.Range("Q3").Value = "$/caja 14 kilos granel"
$ws.Range("S3").Value = 821
$ws.Range("T3").Value = 14

# Row 4: date and volume update
$ws.Range("D4").Value = 44330
$ws.Range("M4").Value = 60

# Row 5: date and volume update
$ws.Range("D5").Value = 44302
$ws.Range("M5").Value = 80

# Row 6: date and volume update
$ws.Range("D6").Value = 44313
$ws.Range("M6").Value = 120

# Row 7: date and volume update
$ws.Range("D7").Value = 44327
$ws.Range("M7").Value = 60

# Row 8: date and volume update
$ws.Range("D8").Value = 44322
$ws.Range("M8").Value = 60

# Row 10: date, volume, and commercialization unit/price details update
$ws.Range("D10").Value = 44316
$ws.Range("M10").Value = 120
$ws.Range("Q10").Value = "$/caja 10 kilos empedrada"
$ws.Range("S10").Value = 11500
$ws.Range("T10").Value = 1
